$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.988.62"
$ws.Range("E2").Value = "  +4.28%  "

$ws.Range("D3").Value = "3.462.01"
$ws.Range("E3").Value = "  +4.22%  "

$ws.Range("D5").Value = "'583.88"
$ws.Range("E5").Value = "  +5.80%  "

$ws.Range("D6").Value = "'186.36"
$ws.Range("E6").Value = "  +7.75%  "

$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("D8").Value = "3.454.87"
$ws.Range("E8").Value = "  +4.30%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").Value = "'0.647"
$ws.Range("E11").Value = "  +2.55%  "

$ws.Range("D12").Value = "'56.04"
$ws.Range("E12").Value = "  +5.75%  "

$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").Value = "'9.41"
$ws.Range("E14").Value = "  +4.09%  "

$ws.Range("D15").Value = "4.025.14"

$ws.Range("D16").Value = "'18.78"
$ws.Range("E16").Value = "  +4.06%  "

$ws.Range("D17").Value = "3.473.20"
$ws.Range("E17").Value = "  +4.53%  "

$ws.Range("D18").Value = "67.075.31"
$ws.Range("E18").Value = "  +4.48%  "

$ws.Range("D19").Value = "'12.15"
$ws.Range("E19").Value = "  +4.24%  "

$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("E21").Value = "  +3.65%  "

$ws.Range("D22").Value = "'485.28"
$ws.Range("E22").Value = "  +7.55%  "

$ws.Range("D23").Value = "'5.31"
$ws.Range("E23").Value = "  +6.37%  "

$ws.Range("D24").Value = "'16.82"
$ws.Range("E24").Value = "  +21.35%  "

$ws.Range("E25").Value = "  +10.83%  "

$ws.Range("D26").Value = "'89.70"
$ws.Range("E26").Value = "  +3.00%  "

$ws.Range("E27").Value = "  +3.29%  "

$ws.Range("D28").Value = "'10.95"
$ws.Range("E28").Value = "  +4.02%  "

$ws.Range("D29").Value = "'9.12"
$ws.Range("E29").Value = "  +6.65%  "

$ws.Range("D30").Value = "'31.35"
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = "  +10.99%  "

$ws.Range("D32").Value = "'600.70"
$ws.Range("E32").Value = "  +5.49%  "

$ws.Range("E33").Value = "  +3.30%  "

$ws.Range("D34").Value = "'63.63"
$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("D35").Value = "'0.111"
$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("D36").Value = "'0.148"
$ws.Range("E36").Value = "  +5.08%  "

$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "'36.47"
$ws.Range("E38").Value = "  +3.80%  "

$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").Value = "'0.384"
$ws.Range("E40").Value = "  +5.49%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.254.30"
$ws.Range("E41").Value = "  +6.29%  "

$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0756"
$ws.Range("E42").Value = "  +4.37%  "

$ws.Range("E43").Value = "  +7.12%  "

$ws.Range("E44").Value = "  +4.37%  "

$ws.Range("D45").Value = "'2.52"
$ws.Range("E45").Value = "  +3.66%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.27"
$ws.Range("E46").Value = "  +2.76%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.74"
$ws.Range("E47").Value = "  +22.59%  "

$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").Value = "'3.28"
$ws.Range("E49").Value = "  +13.31%  "

$ws.Range("D50").Value = "'8.74"
$ws.Range("E50").Value = "  +6.97%  "

$ws.Range("E51").Value = "  +0.16%  "
